$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.843.25'
$ws.Range('D3').Value = '1.857.86'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '304.09'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.5063'
$ws.Range('E7').Value = '  -1.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '0.3646'
$ws.Range('E8').Value = '  -3.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.07164'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = '1.859.51'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '0.07424'
$ws.Range('E13').Value = '  -2.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '92.55'
$ws.Range('E14').Value = '  +3.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '5.225'
$ws.Range('E15').Value = '  -2.23%  '
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '0.000008506'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '14.03'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '0.9995'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').Value = '26.884.56'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('E21').Value = '  -1.22%  '
$ws.Range('D22').Value = '2.093.11'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '10.32'
$ws.Range('E23').Value = '  -3.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '6.423'
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '147.13'
$ws.Range('E25').Value = '  -2.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '1.795'
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '2.056'
$ws.Range('E28').Value = '  -3.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '113.01'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '4.631'
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '4.661'
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '0.09221'
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '0.05078'
$ws.Range('E33').Value = '  -1.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '2.981'
$ws.Range('E34').Value = '  -4.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '0.7443'
$ws.Range('E35').Value = '  -1.42%  '
$ws.Range('E36').Value = '  -2.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '3.242'
$ws.Range('E37').Value = '  +6.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '2.503'
$ws.Range('E38').Value = '  -1.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '0.01990'
$ws.Range('E39').Value = '  -2.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '1.077'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '0.5325'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '118.89'
$ws.Range('E42').Value = '  +3.71%  '
$ws.Range('E43').Value = '  -3.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '8.374'
$ws.Range('E44').Value = '  -2.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '0.1457'
$ws.Range('E45').Value = '  -1.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '0.4642'
$ws.Range('E46').Value = '  -0.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '0.9994'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '9.996'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '1.559'
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '36.78'
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '62.85'
$ws.Range('E51').Value = '  -3.65%  '
